# The deck's two theme parts (ppt/theme/theme1.xml = "Office Theme",
# ppt/theme/theme2.xml = "Integral") had their contents swapped: theme2.xml
# (the theme actually wired to the slide master / presentation) becomes the
# plain "Office Theme" palette, while theme1.xml becomes "Integral". The
# PowerPoint object model only exposes a single live Theme object (the one
# backing the SlideMaster/Presentation, i.e. ppt/theme/theme2.xml), so we
# repoint that object's color scheme (and its name) at the "Office Theme"
# values described by the diff.

$p = $ppt.ActivePresentation
$theme = $p.SlideMaster.Theme
$tcs = $theme.ThemeColorScheme

# Index order exposed by ThemeColorScheme.Colors: 1 dk1, 2 lt1, 3 dk2, 4 lt2,
# 5-10 accent1-6, 11 hlink, 12 folHlink.
$tcs.Colors(1).RGB  = 0        # dk1      000000
$tcs.Colors(2).RGB  = 16777215 # lt1      FFFFFF
$tcs.Colors(3).RGB  = 6968388  # dk2      44546A
$tcs.Colors(4).RGB  = 15132391 # lt2      E7E6E6
$tcs.Colors(5).RGB  = 13998939 # accent1  5B9BD5
$tcs.Colors(6).RGB  = 3243501  # accent2  ED7D31
$tcs.Colors(7).RGB  = 10855845 # accent3  A5A5A5
$tcs.Colors(8).RGB  = 49407    # accent4  FFC000
$tcs.Colors(9).RGB  = 12874308 # accent5  4472C4
$tcs.Colors(10).RGB = 4697456  # accent6  70AD47
$tcs.Colors(11).RGB = 12673797 # hlink    0563C1
$tcs.Colors(12).RGB = 7491477  # folHlink 954F72

$theme.Name = "Office Theme"
